# Auto-generated edit script
# Applies explicit cell value updates (and clears) per the target diff,
# derived from a scheduled market-data refresh across the Phoenix_Profits sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null
$ws.Range("H32").Value = 2697
$ws.Range("I32").Value = 2020.5
$ws.Range("J32").Value = 3869.6
$ws.Range("K32").Value = 2020.5
$ws.Range("L32").Value = 3869.6
$ws.Range("M32").Value = -1694.5
$ws.Range("N32").Value = -4521.6
$ws.Range("H80").Value = 761.95
$ws.Range("I80").Value = 934.8182
$ws.Range("J80").Value = 550.6667
$ws.Range("K80").Value = 2804.4546
$ws.Range("L80").Value = 1652.0001
$ws.Range("M80").Value = -1806.4546
$ws.Range("N80").Value = -3648.0001
$ws.Range("H83").Value = 761.95
$ws.Range("I83").Value = 934.8182
$ws.Range("J83").Value = 550.6667
$ws.Range("K83").Value = 8413.363800000001
$ws.Range("L83").Value = 4956.0003
$ws.Range("M83").Value = -3421.363800000001
$ws.Range("N83").Value = -14940.0003
$ws.Range("H95").Value = 54899.1
$ws.Range("J95").Value = 54899.1
$ws.Range("L95").Value = 54899.1
$ws.Range("N95").Value = -60391.1
$ws.Range("H112").Value = 84012.91
$ws.Range("I112").Value = 1123.5
$ws.Range("J112").Value = 102432.78
$ws.Range("K112").Value = 3370.5
$ws.Range("L112").Value = 307298.34
$ws.Range("M112").Value = -2262.5
$ws.Range("N112").Value = -309514.34
$ws.Range("H113").Value = 2130.111
$ws.Range("J113").Value = 2239.2
$ws.Range("L113").Value = 2239.2
$ws.Range("N113").Value = -8747.200000000001
$ws.Range("H116").Value = 7801.0557
$ws.Range("I116").Value = 7329.3335
$ws.Range("J116").Value = 8036.9165
$ws.Range("K116").Value = 7329.3335
$ws.Range("L116").Value = 8036.9165
$ws.Range("M116").Value = -3887.3335
$ws.Range("N116").Value = -14920.9165
$ws.Range("H129").Value = 28728.771
$ws.Range("I129").Value = 65598.42999999999
$ws.Range("J129").Value = 4149
$ws.Range("K129").Value = 196795.29
$ws.Range("L129").Value = 12447
$ws.Range("M129").Value = -191795.29
$ws.Range("N129").Value = -22447
$ws.Range("H132").Value = 2251.372
$ws.Range("I132").Value = 1758.1316
$ws.Range("K132").Value = 5274.3948
$ws.Range("M132").Value = -2744.3948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = $null
$ws.Range("H32").Value = 19673.57
$ws.Range("I32").Value = 21465.037
$ws.Range("K32").Value = 21465.037
$ws.Range("M32").Value = -21178.037
$ws.Range("H51").Value = 35165.668
$ws.Range("J51").Value = 35165.668
$ws.Range("L51").Value = 35165.668
$ws.Range("N51").Value = -36677.668
$ws.Range("H61").Value = 39536.85
$ws.Range("I61").Value = 4508
$ws.Range("K61").Value = 4508
$ws.Range("M61").Value = -4296
$ws.Range("H132").Value = 3906.2222
$ws.Range("I132").Value = 3906.2222
$ws.Range("K132").Value = 11718.6666
$ws.Range("M132").Value = -9188.6666
$ws.Range("H136").Value = 39536.85
$ws.Range("I136").Value = 4508
$ws.Range("K136").Value = 13524
$ws.Range("M136").Value = -10974

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 382.33334
$ws.Range("I22").Value = 501
$ws.Range("K22").Value = 501
$ws.Range("M22").Value = -328
$ws.Range("H105").Value = 2225.35
$ws.Range("I105").Value = 1883.75
$ws.Range("J105").Value = 3591.75
$ws.Range("K105").Value = 1883.75
$ws.Range("L105").Value = 3591.75
$ws.Range("M105").Value = -136.75
$ws.Range("N105").Value = -7085.75
$ws.Range("H134").Value = 32038.291
$ws.Range("I134").Value = 39933.875
$ws.Range("K134").Value = 119801.625
$ws.Range("M134").Value = -117266.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 25000
$ws.Range("I13").Value = 25000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 25000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -24861
$ws.Range("N13").Value = $null
$ws.Range("H31").Value = 2376.7593
$ws.Range("J31").Value = 3561.5334
$ws.Range("L31").Value = 3561.5334
$ws.Range("N31").Value = -4151.5334
$ws.Range("H34").Value = 2376.7593
$ws.Range("J34").Value = 3561.5334
$ws.Range("L34").Value = 3561.5334
$ws.Range("N34").Value = -3965.5334
$ws.Range("H68").Value = 38667.69
$ws.Range("J68").Value = 38556.668
$ws.Range("L68").Value = 38556.668
$ws.Range("N68").Value = -40054.668
$ws.Range("H71").Value = 38667.69
$ws.Range("J71").Value = 38556.668
$ws.Range("L71").Value = 115670.004
$ws.Range("N71").Value = -123158.004
$ws.Range("H99").Value = 6415.05
$ws.Range("I99").Value = 5277.75
$ws.Range("J99").Value = 8121
$ws.Range("K99").Value = 5277.75
$ws.Range("L99").Value = 8121
$ws.Range("M99").Value = -3779.75
$ws.Range("N99").Value = -11117
$ws.Range("H126").Value = 6415.05
$ws.Range("I126").Value = 5277.75
$ws.Range("J126").Value = 8121
$ws.Range("K126").Value = 15833.25
$ws.Range("L126").Value = 24363
$ws.Range("M126").Value = -13363.25
$ws.Range("N126").Value = -29303
$ws.Range("H134").Value = 4206.227
$ws.Range("I134").Value = 4145.3687
$ws.Range("K134").Value = 12436.1061
$ws.Range("M134").Value = -9901.106100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 633.1539
$ws.Range("I8").Value = 633.1539
$ws.Range("K8").Value = 1899.4617
$ws.Range("M8").Value = -1760.4617
$ws.Range("H63").Value = 8000
$ws.Range("I63").Value = 8000
$ws.Range("K63").Value = 24000
$ws.Range("M63").Value = -23251
$ws.Range("H66").Value = 8000
$ws.Range("I66").Value = 8000
$ws.Range("K66").Value = 72000
$ws.Range("M66").Value = -68256
$ws.Range("H81").Value = 1996.25
$ws.Range("I81").Value = 1996.25
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5988.75
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4865.75
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 1996.25
$ws.Range("I84").Value = 1996.25
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 17966.25
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -12350.25
$ws.Range("N84").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 32459
$ws.Range("J15").Value = 32459
$ws.Range("L15").Value = 32459
$ws.Range("N15").Value = -33035
$ws.Range("H81").Value = 32459
$ws.Range("J81").Value = 32459
$ws.Range("L81").Value = 32459
$ws.Range("N81").Value = -34455
$ws.Range("H84").Value = 32459
$ws.Range("J84").Value = 32459
$ws.Range("L84").Value = 97377
$ws.Range("N84").Value = -107361
$ws.Range("H122").Value = 3327.2334
$ws.Range("I122").Value = 2779.4666
$ws.Range("K122").Value = 8338.399800000001
$ws.Range("M122").Value = -5888.399800000001
$ws.Range("H132").Value = 2005374.8
$ws.Range("I132").Value = 2358176.8
$ws.Range("K132").Value = 7074530.399999999
$ws.Range("M132").Value = -7072000.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 39971.29
$ws.Range("I40").Value = 45531.25
$ws.Range("J40").Value = 11059.5
$ws.Range("K40").Value = 45531.25
$ws.Range("L40").Value = 11059.5
$ws.Range("M40").Value = -45395.25
$ws.Range("N40").Value = -11331.5
$ws.Range("H46").Value = 4085.1667
$ws.Range("I46").Value = 771.6667
$ws.Range("J46").Value = 7398.6665
$ws.Range("K46").Value = 771.6667
$ws.Range("L46").Value = 7398.6665
$ws.Range("M46").Value = -583.6667
$ws.Range("N46").Value = -7774.6665
$ws.Range("H61").Value = 3221.3333
$ws.Range("I61").Value = 3221.3333
$ws.Range("K61").Value = 3221.3333
$ws.Range("M61").Value = -3019.3333
$ws.Range("H93").Value = 2854.7856
$ws.Range("I93").Value = 2706.9167
$ws.Range("J93").Value = 3742
$ws.Range("K93").Value = 2706.9167
$ws.Range("L93").Value = 3742
$ws.Range("M93").Value = -1458.9167
$ws.Range("N93").Value = -6238
$ws.Range("H113").Value = 3221.3333
$ws.Range("I113").Value = 3221.3333
$ws.Range("K113").Value = 3221.3333
$ws.Range("M113").Value = -1051.3333
$ws.Range("H132").Value = 3219.1956
$ws.Range("I132").Value = 2866.75
$ws.Range("K132").Value = 8600.25
$ws.Range("M132").Value = -6070.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1666.1666
$ws.Range("I81").Value = 1031.1538
$ws.Range("J81").Value = 3317.2
$ws.Range("K81").Value = 2062.3076
$ws.Range("L81").Value = 6634.4
$ws.Range("M81").Value = -1001.3076
$ws.Range("N81").Value = -8756.4
$ws.Range("H84").Value = 1666.1666
$ws.Range("I84").Value = 1031.1538
$ws.Range("J84").Value = 3317.2
$ws.Range("K84").Value = 10311.538
$ws.Range("L84").Value = 33172
$ws.Range("M84").Value = -5007.538
$ws.Range("N84").Value = -43780
$ws.Range("H95").Value = 74999.336
$ws.Range("J95").Value = 74999.336
$ws.Range("L95").Value = 74999.336
$ws.Range("N95").Value = -80491.336
$ws.Range("H132").Value = 29586.2
$ws.Range("I132").Value = 29566.285
$ws.Range("J132").Value = 29632.666
$ws.Range("K132").Value = 88698.855
$ws.Range("L132").Value = 88897.99800000001
$ws.Range("M132").Value = -86168.855
$ws.Range("N132").Value = -93957.99800000001
$ws.Range("H133").Value = 79999
$ws.Range("J133").Value = 79999
$ws.Range("L133").Value = 79999
$ws.Range("N133").Value = -90119
$ws.Range("H136").Value = 3537.6785
$ws.Range("I136").Value = 3474.28
$ws.Range("K136").Value = 10422.84
$ws.Range("M136").Value = -7872.84
